# Update Name of Algo
# Apply updated values to specific cells in columns A and B as per the
# refreshed KNN imputation result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 7.733999999999999
$ws.Range("A9").Value  = -21.462
$ws.Range("B9").Value  = 6.67
$ws.Range("A18").Value = -21.995
$ws.Range("A20").Value = -20.511
$ws.Range("B23").Value = 7.62
$ws.Range("B24").Value = 5.453
$ws.Range("B26").Value = 5.559
$ws.Range("A27").Value = -21.791
$ws.Range("B34").Value = 7.447
$ws.Range("B35").Value = 8.164
$ws.Range("B48").Value = 5.616
$ws.Range("B52").Value = 5.573
$ws.Range("B66").Value = 5.202
$ws.Range("B67").Value = 5.329
$ws.Range("A69").Value = -21.452
$ws.Range("A76").Value = -20.468
$ws.Range("B80").Value = 8.687999999999999
$ws.Range("A82").Value = -21.81
$ws.Range("B99").Value = 5.206999999999999
